$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 3463
$ws.Range("F4").Value = 368
$ws.Range("F5").Value = 8109
$ws.Range("F7").Value = 66
$ws.Range("F9").Value = 65
$ws.Range("F14").Value = 1061
$ws.Range("F18").Value = 1139
$ws.Range("F20").Value = 717
$ws.Range("F21").Value = 509
$ws.Range("F24").Value = 4634
$ws.Range("F25").Value = 102
$ws.Range("C26").Value = '上海·【开票倒计时】Wonder Festival 2024上海'
$ws.Range("F26").Value = 49281
$ws.Range("F27").Value = 3989
$ws.Range("F30").Value = 739
$ws.Range("F31").Value = 144
$ws.Range("F35").Value = 565
$ws.Range("F36").Value = 189
$ws.Range("F40").Value = 986
$ws.Range("F42").Value = 153
$ws.Range("F43").Value = 1047
$ws.Range("F44").Value = 675
$ws.Range("F45").Value = 90
$ws.Range("F47").Value = 84
$ws.Range("F48").Value = 22
$ws.Range("F49").Value = 2451
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 239
$ws.Range("F18").Value = 25
$ws.Range("F19").Value = 142
$ws.Range("F20").Value = 7284
$ws.Range("F30").Value = 74
$ws = $wb.Worksheets.Item(3)
$ws.Range("F5").Value = 1476
$ws.Range("F8").Value = 2315
$ws.Range("F9").Value = 9238
$ws.Range("F10").Value = 1506
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 3463
$ws.Range("F4").Value = 8109
$ws.Range("F5").Value = 1476
$ws.Range("F7").Value = 1506
$ws.Range("F9").Value = 66
$ws.Range("F10").Value = 65
$ws.Range("F13").Value = 1061
$ws.Range("F14").Value = 239
$ws.Range("F19").Value = 1139
$ws.Range("F22").Value = 4634
$ws.Range("F23").Value = 102
$ws.Range("F26").Value = 3989
$ws.Range("F29").Value = 739
$ws.Range("F30").Value = 144
$ws.Range("F33").Value = 565
$ws.Range("F35").Value = 189
$ws.Range("F40").Value = 153
$ws.Range("F41").Value = 1047
$ws.Range("F42").Value = 675
$ws.Range("F44").Value = 90
$ws.Range("F46").Value = 84
$ws.Range("F47").Value = 22
